$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '76.562.04'
$ws.Range('E2').Value = '  +0.64%  '

$ws.Range('D3').Value = '3.039.29'
$ws.Range('E3').Value = '  +4.20%  '

$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('D5').Value = "'202.09"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.30%  '

$ws.Range('D6').Value = "'629.83"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +5.33%  '

$ws.Range('E7').Value = '  +0.06%  '

$ws.Range('D8').Value = "'0.553"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +0.64%  '

$ws.Range('D9').Value = "'0.210"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +6.11%  '

$ws.Range('D10').Value = '3.038.62'
$ws.Range('E10').Value = '  +4.28%  '

$ws.Range('E11').Value = '  +2.07%  '

$ws.Range('E12').Value = '  -0.49%  '

$ws.Range('D13').Value = "'5.12"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +5.25%  '

$ws.Range('D14').Value = '3.600.66'

$ws.Range('D15').Value = "'29.52"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +6.57%  '

$ws.Range('D16').Value = '76.481.33'
$ws.Range('E16').Value = '  +0.72%  '

$ws.Range('E17').Value = '  +2.03%  '

$ws.Range('D18').Value = '3.031.78'
$ws.Range('E18').Value = '  +4.09%  '

$ws.Range('D19').Value = "'13.46"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +4.24%  '

$ws.Range('D20').Value = "'9.06"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.60%  '

$ws.Range('D21').Value = "'376.55"
$ws.Range('D21').Style = 'Normal'

$ws.Range('E22').Value = '  -0.36%  '

$ws.Range('E23').Value = '  +2.14%  '

$ws.Range('D24').Value = "'73.81"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.06%  '

$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').Value = "'1.00"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.19%  '

$ws.Range('B27').Value = 'NEARProtocol'
$ws.Range('C27').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D27').Value = "'4.39"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.94%  '

$ws.Range('D28').Value = "'10.00"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.18%  '

$ws.Range('E29').Value = '  +4.11%  '

$ws.Range('E30').Value = '  +0.25%  '

$ws.Range('E31').Value = '  +8.28%  '

$ws.Range('E32').Value = '  +1.32%  '

$ws.Range('D33').Value = "'515.97"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.74%  '

$ws.Range('D34').Value = "'1.96"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +8.28%  '

$ws.Range('D35').Value = "'1.00"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.04%  '

$ws.Range('D36').Value = "'20.86"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +3.35%  '

$ws.Range('D37').Value = "'163.53"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.95%  '

$ws.Range('E38').Value = '  +10.43%  '

$ws.Range('D39').Value = "'20.03"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +2.09%  '

$ws.Range('E40').Value = '  +3.24%  '

$ws.Range('D41').Value = "'188.29"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.91%  '

$ws.Range('D42').Value = "'0.112"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.48%  '

$ws.Range('E43').Value = '  +0.01%  '

$ws.Range('E44').Value = '  +4.76%  '

$ws.Range('D45').Value = "'1.27"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +6.03%  '

$ws.Range('D46').Value = "'42.11"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +4.75%  '

$ws.Range('D47').Value = "'1.67"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.16%  '

$ws.Range('D48').Value = "'0.730"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +11.13%  '

$ws.Range('D49').Value = "'2.44"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +4.08%  '

$ws.Range('D50').Value = "'0.607"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +6.34%  '

$ws.Range('E51').Value = '  +4.62%  '
